$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "456/CASA"
$ws.Range("B2").Value = "Direction régionale"
$ws.Range("C2").Value = "BK747A53"
$ws.Range("D2").Value = "ANAS MASTI"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("N2").Value = "--"
$ws.Range("O2").Value = 1000

# Row 3
$ws.Range("A3").Value = "456/CASA"
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "BG746583"
$ws.Range("D3").Value = "TEST KAMAL"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 1000
$ws.Range("J3").Value = 0
$ws.Range("O3").Value = 1000

# Row 4
$ws.Range("A4").Value = "456/CASA"
$ws.Range("C4").Value = "GT744635"
$ws.Range("D4").Value = "JAMAL JAMAL"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1000
$ws.Range("J4").Value = 0
$ws.Range("O4").Value = 1000

# Row 5
$ws.Range("H5").Value = 3000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 3000
